# Merge 'TSM Alistair's Rose'
# Adds a new strings-table entry (item name + description) for the new
# "Alistair's Rose" amulet, with a matching developer comment on the new
# name cell, and nudges the saved window/selection state the way the
# authoring session left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows -----------------------------------------------------
# Row 58: string id + the item's display name ("Alistair's Rose")
# Row 59: string id + the item's (wrapped, taller) flavour-text description
$ws.Range("A58").Value = 6610056
$ws.Range("B58").Value = "Alistair's Rose"

$ws.Range("A59").Value = 6610057
$ws.Range("B59").Value = "A beautiful rose with velvety petals of deep red. This was a gift from Alistair.`nWhether infused with lyrium, magically preserved, or just the perfect cultivar for cut flower arrangements, the enduring loveliness of this single perfect bloom serves as a reminder that rare and wonderful things can indeed be found amidst all the darkness."

# --- Formatting: reuse the existing "new content" look ------------------
# Column A on these rows uses the same highlighted style as rows 3-5
# (the most recently added id column cells).
$ws.Range("A3").Copy()
$ws.Range("A58:A59").PasteSpecial(-4122)

# B58 is a plain single-line entry like the rest of the name column.
$ws.Range("B3").Copy()
$ws.Range("B58").PasteSpecial(-4122)

# B59 is a wrapped, multi-line description like the other item descriptions.
$ws.Range("B13").Copy()
$ws.Range("B59").PasteSpecial(-4122)

# Give the description row the taller height used for wrapped description rows.
$ws.Rows.Item(59).RowHeight = 45

# --- Reviewer comment on the new name cell -------------------------------
$ws.Range("A58").AddComment("TSM Alistairs Rose") | Out-Null

# --- Restore window/selection state left by the authoring session -------
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Left = 30075
$excel.ActiveWindow.Top = 2010
$ws.Range("B73").Select() | Out-Null
